$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.179.10'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '3.141.23'
$ws.Range('E3').Value = '  +2.38%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.13'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '180.41'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +7.18%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.138.42'
$ws.Range('E8').Value = '  +2.36%  '
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.56'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.153'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.50%  '
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.93'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.94%  '
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '3.665.35'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = '68.086.54'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.14'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('D19').Value = '3.137.86'
$ws.Range('E19').Value = '  +2.38%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.55'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '490.30'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  +1.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.81'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '84.05'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.36'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +7.80%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.97'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.70'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.54%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.18'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +5.51%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.36'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.63%  '
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '28.38'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.36%  '
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('D34').Value = '0.0₃0958'
$ws.Range('E34').Value = '  +5.53%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '48.88'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +5.11%  '
$ws.Range('E37').Value = '  +2.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.962'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.325'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +8.68%  '
$ws.Range('E40').Value = '  +4.87%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '49.37'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.125'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.43'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.71'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +9.81%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '393.51'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +7.18%  '
$ws.Range('D46').Value = '2.789.34'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.99'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +10.67%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0349'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '136.71'
$ws.Range('D49').Style = "Normal"
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.34'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +8.89%  '
